$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (Changed) date column (C) for all data rows (2-52)
# from serial date 45190 to 45192.
$ws.Range("C2:C52").Value = 45192
